# Update row 2 results for each year sheet with new server values.
$wb = $excel.ActiveWorkbook

$data = @{
    "2025" = @{ "B" = 277.8902526399997; "E" = 29092.72506141524; "G" = 8095.925712664175; "I" = 14865.25872276; "L" = 50912.59821312752; "M" = 11247.09127927; "N" = 7270.40043739626; "O" = 6890.515476413165 }
    "2030" = @{ "A" = 0; "B" = 5522.228665637358; "E" = 56033.89983931746; "G" = 8095.925712664175; "I" = 37494.07304221769; "L" = 72580.41557995854; "M" = 21817.92954626511; "N" = 10913.45045002708; "O" = 9426.543084076415 }
    "2035" = @{ "A" = 2266.487971660625; "B" = 7059.768887134273; "E" = 67081.11907838765; "G" = 8095.925712664175; "I" = 53779.76068266826; "L" = 72580.41557995854; "M" = 27679.17349794621; "N" = 15939.31470367384; "O" = 15294.37937820872 }
    "2040" = @{ "A" = 2266.487971660625; "B" = 7059.768887134273; "E" = 67081.11907838765; "G" = 8095.925712664175; "I" = 53779.76068266826; "L" = 72580.41557995854; "M" = 27679.17349794621; "N" = 15939.31470367384; "O" = 15294.37937820872 }
    "2045" = @{ "A" = 2266.487971660625; "B" = 7059.768887134273; "E" = 67081.11907838765; "G" = 8095.925712664175; "I" = 53779.76068266826; "L" = 72580.41557995854; "M" = 27679.17349794621; "N" = 15939.31470367384; "O" = 15294.37937820872 }
    "2050" = @{ "A" = 2266.487971660625; "B" = 7059.768887134273; "E" = 67081.11907838765; "G" = 8095.925712664175; "I" = 53779.76068266826; "L" = 72580.41557995854; "M" = 27679.17349794621; "N" = 15939.31470367384; "O" = 15294.37937820872 }
}

foreach ($sheetKey in $data.Keys) {
    $sheetName = [string]$sheetKey
    $ws = $wb.Worksheets.Item($sheetName)
    $cols = $data[$sheetKey]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col`2").Value = $cols[$col]
    }
}
